# Update countries & provincias Spain
# Refresh "Pais" COVID-19 stats sheet with the latest scraped totals and
# refresh the "last updated" timestamp. A handful of countries (Hungria,
# Bulgaria, Consejo Danes para los Refugiados, Martinica) overtook their
# neighbour in the Casos totales ranking, so their row now carries the
# other country's name/data and vice versa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "Datos actualizados..." banner -> refreshed timestamp (07:22 -> 07:52)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 24 de Abril de 2020 a las 07:52"

# Row 46: Australia
$ws.Cells.Item(46, 4).Value = 5129
$ws.Cells.Item(46, 5).Value = 1468

# Row 48: Republica Dominicana
$ws.Cells.Item(48, 4).Value = 674
$ws.Cells.Item(48, 5).Value = 4604

# Row 61: Grecia
$ws.Cells.Item(61, 5).Value = 1759
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 127

# Row 63: Hungria
$ws.Cells.Item(63, 1).Value = "Hungria"
$ws.Cells.Item(63, 2).Value = 2383
$ws.Cells.Item(63, 3).Value = 99
$ws.Cells.Item(63, 4).Value = 401
$ws.Cells.Item(63, 5).Value = 1732
$ws.Cells.Item(63, 6).Value = 61
$ws.Cells.Item(63, 7).Value = 11
$ws.Cells.Item(63, 8).Value = 250

# Row 64: Kazajistan
$ws.Cells.Item(64, 1).Value = "Kazajistan"
$ws.Cells.Item(64, 2).Value = 2334
$ws.Cells.Item(64, 3).Value = 45
$ws.Cells.Item(64, 4).Value = 560
$ws.Cells.Item(64, 5).Value = 1752
$ws.Cells.Item(64, 6).Value = 29
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = 22

# Row 68: Uzbekistan
$ws.Cells.Item(68, 2).Value = 1778
$ws.Cells.Item(68, 3).Value = 20
$ws.Cells.Item(68, 4).Value = 563
$ws.Cells.Item(68, 5).Value = 1208

# Row 83: Bulgaria
$ws.Cells.Item(83, 1).Value = "Bulgaria"
$ws.Cells.Item(83, 2).Value = 1171
$ws.Cells.Item(83, 3).Value = 74
$ws.Cells.Item(83, 4).Value = 193
$ws.Cells.Item(83, 5).Value = 926
$ws.Cells.Item(83, 6).Value = 37
$ws.Cells.Item(83, 8).Value = 52

# Row 84: Ghana
$ws.Cells.Item(84, 1).Value = "Ghana"
$ws.Cells.Item(84, 2).Value = 1154
$ws.Cells.Item(84, 4).Value = 99
$ws.Cells.Item(84, 5).Value = 1046
$ws.Cells.Item(84, 6).Value = 4
$ws.Cells.Item(84, 8).Value = 9

# Row 102: Honduras
$ws.Cells.Item(102, 4).Value = 50
$ws.Cells.Item(102, 5).Value = 465

# Row 112: Consejo Danes para los Refugiados
$ws.Cells.Item(112, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(112, 2).Value = 394
$ws.Cells.Item(112, 3).Value = 17
$ws.Cells.Item(112, 4).Value = 48
$ws.Cells.Item(112, 5).Value = 321
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 25

# Row 113: Guatemala
$ws.Cells.Item(113, 1).Value = "Guatemala"
$ws.Cells.Item(113, 2).Value = 384
$ws.Cells.Item(113, 3).Value = 42
$ws.Cells.Item(113, 4).Value = 30
$ws.Cells.Item(113, 5).Value = 343
$ws.Cells.Item(113, 6).Value = 5
$ws.Cells.Item(113, 7).Value = 1
$ws.Cells.Item(113, 8).Value = 11

# Row 124: Vietnam
$ws.Cells.Item(124, 4).Value = 225
$ws.Cells.Item(124, 5).Value = 43

# Row 131: Martinica
$ws.Cells.Item(131, 1).Value = "Martinica"
$ws.Cells.Item(131, 2).Value = 170
$ws.Cells.Item(131, 3).Value = 6
$ws.Cells.Item(131, 4).Value = 77
$ws.Cells.Item(131, 5).Value = 79
$ws.Cells.Item(131, 6).Value = 6
$ws.Cells.Item(131, 8).Value = 14

# Row 132: Gabon
$ws.Cells.Item(132, 1).Value = "Gabon"
$ws.Cells.Item(132, 2).Value = 167
$ws.Cells.Item(132, 4).Value = 24
$ws.Cells.Item(132, 5).Value = 141
$ws.Cells.Item(132, 6).Value = 1
$ws.Cells.Item(132, 8).Value = 2

# Row 134: Guadalupe
$ws.Cells.Item(134, 2).Value = 149
$ws.Cells.Item(134, 3).Value = 1
$ws.Cells.Item(134, 4).Value = 82
$ws.Cells.Item(134, 5).Value = 55
$ws.Cells.Item(134, 6).Value = 11

# Row 158: Guyana
$ws.Cells.Item(158, 4).Value = 12
$ws.Cells.Item(158, 5).Value = 51

# Row 166: Mozambique
$ws.Cells.Item(166, 4).Value = 12
$ws.Cells.Item(166, 5).Value = 34

# Row 191: Granada
$ws.Cells.Item(191, 4).Value = 7
$ws.Cells.Item(191, 5).Value = 8
